$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 80, pushing the existing rows
# (old 80-165) down to 82-167, matching the new dimension A1:R167.
$ws.Rows("80:81").Insert()

# Populate the newly inserted row 80 with its data.
$ws.Cells.Item(80, 1).Value = 10
$ws.Cells.Item(80, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(80, 3).Value = "La Araucanía"
$ws.Cells.Item(80, 4).Value = 44494
$ws.Cells.Item(80, 5).Value = 9
$ws.Cells.Item(80, 6).Value = 100112039
$ws.Cells.Item(80, 7).Value = "Ciboulette"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 20
$ws.Cells.Item(80, 11).Value = 6000
$ws.Cells.Item(80, 12).Value = 6000
$ws.Cells.Item(80, 13).Value = 6000
$ws.Cells.Item(80, 14).Value = '$/docena de atados'
$ws.Cells.Item(80, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(80, 16).Value = 2000
$ws.Cells.Item(80, 17).Value = 3
$ws.Cells.Item(80, 18).Value = "Hortaliza"

# Populate the newly inserted row 81 with its data.
$ws.Cells.Item(81, 1).Value = 10
$ws.Cells.Item(81, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value = "La Araucanía"
$ws.Cells.Item(81, 4).Value = 44494
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = 100112039
$ws.Cells.Item(81, 7).Value = "Ciboulette"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 30
$ws.Cells.Item(81, 11).Value = 2000
$ws.Cells.Item(81, 12).Value = 2000
$ws.Cells.Item(81, 13).Value = 2000
$ws.Cells.Item(81, 14).Value = '$/docena de atados'
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 667
$ws.Cells.Item(81, 17).Value = 3
$ws.Cells.Item(81, 18).Value = "Hortaliza"
